$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.58242506653322
$ws.Range("C2").Value = 11.17072946123901
$ws.Range("D2").Value = 6.035509846248532
$ws.Range("E2").Value = 9.992111519605857
$ws.Range("G2").Value = 3.715378354171635
$ws.Range("I2").Value = 31.20382436938653
$ws.Range("K2").Value = 14.77870675430539
$ws.Range("M2").Value = 16.44242788566918
$ws.Range("N2").Value = 23.24136519306006

$ws.Range("B3").Value = 13.35983330404993
$ws.Range("C3").Value = 10.92935357479295
$ws.Range("D3").Value = 5.926274958308508
$ws.Range("E3").Value = 9.784645733394083
$ws.Range("G3").Value = 3.718959085805852
$ws.Range("I3").Value = 31.089795164705
$ws.Range("K3").Value = 14.61237133735311
$ws.Range("M3").Value = 16.29035581603003
$ws.Range("N3").Value = 23.23816691146897

$ws.Range("B4").Value = 13.22583735744279
$ws.Range("C4").Value = 10.7820280801196
$ws.Range("D4").Value = 5.860112251089164
$ws.Range("E4").Value = 9.658427016135146
$ws.Range("G4").Value = 3.721270185053882
$ws.Range("I4").Value = 31.02468429041672
$ws.Range("K4").Value = 14.5137719216285
$ws.Range("M4").Value = 16.20108916562506
$ws.Range("N4").Value = 23.23788174295171

$ws.Range("B5").Value = 13.17198345065281
$ws.Range("C5").Value = 10.72230469026524
$ws.Range("D5").Value = 5.833417559173391
$ws.Range("E5").Value = 9.607362030171597
$ws.Range("G5").Value = 3.72224037930883
$ws.Range("I5").Value = 30.999395496709
$ws.Range("K5").Value = 14.4745251740185
$ws.Range("M5").Value = 16.16578042985412
$ws.Range("N5").Value = 23.23818727621689

$ws.Range("B6").Value = 13.16308862294391
$ws.Range("C6").Value = 10.71240926759761
$ws.Range("D6").Value = 5.829002192160041
$ws.Range("E6").Value = 9.598907341773497
$ws.Range("G6").Value = 3.722403197986848
$ws.Range("I6").Value = 30.99527177725687
$ws.Range("K6").Value = 14.46806593546523
$ws.Range("M6").Value = 16.15998294970997
$ws.Range("N6").Value = 23.23826347014861

$ws.Range("B7").Value = 13.22510792442555
$ws.Range("C7").Value = 10.78122124058966
$ws.Range("D7").Value = 5.859751104257461
$ws.Range("E7").Value = 9.657736735069362
$ws.Range("G7").Value = 3.721283154288592
$ws.Range("I7").Value = 31.02433818445372
$ws.Range("K7").Value = 14.51323878898742
$ws.Range("M7").Value = 16.20060860948099
$ws.Range("N7").Value = 23.23788415648728

$ws.Range("B8").Value = 13.50516680330642
$ws.Range("C8").Value = 11.08737486396151
$ws.Range("D8").Value = 5.997680537203114
$ws.Range("E8").Value = 9.920381249492511
$ws.Range("G8").Value = 3.716589703963101
$ws.Range("I8").Value = 31.16349659868983
$ws.Range("K8").Value = 14.72064924004343
$ws.Range("M8").Value = 16.38916324400694
$ws.Range("N8").Value = 23.23991383635909

$ws.Range("B9").Value = 14.07215244762709
$ws.Range("C9").Value = 11.69079095301951
$ws.Range("D9").Value = 6.273660504996846
$ws.Range("E9").Value = 10.44134639904463
$ws.Range("G9").Value = 3.708273618556911
$ws.Range("I9").Value = 31.47473803112514
$ws.Range("K9").Value = 15.1533128585997
$ws.Range("M9").Value = 16.78982985850837
$ws.Range("N9").Value = 23.25722896018469

$ws.Range("B10").Value = 14.49495908945291
$ws.Range("C10").Value = 12.1309010424839
$ws.Range("D10").Value = 6.477528467507958
$ws.Range("E10").Value = 10.82337164632979
$ws.Range("G10").Value = 3.702697980586468
$ws.Range("I10").Value = 31.72596768662151
$ws.Range("K10").Value = 15.48410547554685
$ws.Range("M10").Value = 17.10068152782674
$ws.Range("N10").Value = 23.27810216554676

$ws.Range("B11").Value = 14.68770575531689
$ws.Range("C11").Value = 12.32941563050344
$ws.Range("D11").Value = 6.570055200391582
$ws.Range("E11").Value = 10.99614034561774
$ws.Range("G11").Value = 3.700275973756063
$ws.Range("I11").Value = 31.84495709523587
$ws.Range("K11").Value = 15.63676306509198
$ws.Range("M11").Value = 17.24514370095895
$ws.Range("N11").Value = 23.28937154302067

$ws.Range("B12").Value = 14.76067600407725
$ws.Range("C12").Value = 12.40426664398811
$ws.Range("D12").Value = 6.605025538141324
$ws.Range("E12").Value = 11.06134951490199
$ws.Range("G12").Value = 3.699375155080495
$ws.Range("I12").Value = 31.89067273564956
$ws.Range("K12").Value = 15.69483045016303
$ws.Range("M12").Value = 17.30024110637907
$ws.Range("N12").Value = 23.29389412181357

$ws.Range("B13").Value = 14.74496259506959
$ws.Range("C13").Value = 12.38816166671574
$ws.Range("D13").Value = 6.597497621511731
$ws.Range("E13").Value = 11.04731614563717
$ws.Range("G13").Value = 3.699568437316376
$ws.Range("I13").Value = 31.88079815851081
$ws.Range("K13").Value = 15.68231392451523
$ws.Range("M13").Value = 17.2883581723601
$ws.Range("N13").Value = 23.29290875825188

$ws.Range("B14").Value = 14.69370987187053
$ws.Range("C14").Value = 12.3355805578378
$ws.Range("D14").Value = 6.57293378927203
$ws.Range("E14").Value = 11.00150982203092
$ws.Range("G14").Value = 3.700201535927572
$ws.Range("I14").Value = 31.84870509787449
$ws.Range("K14").Value = 15.64153541565759
$ws.Range("M14").Value = 17.24966899470029
$ws.Range("N14").Value = 23.28973850579132

$ws.Range("B15").Value = 14.66231140234343
$ws.Range("C15").Value = 12.30332890641283
$ws.Range("D15").Value = 6.557877875174434
$ws.Range("E15").Value = 10.97342221576212
$ws.Range("G15").Value = 3.700591452310842
$ws.Range("I15").Value = 31.82913212762606
$ws.Range("K15").Value = 15.61658956947484
$ws.Range("M15").Value = 17.22602051295473
$ws.Range("N15").Value = 23.28782985617356

$ws.Range("B16").Value = 14.48236418724899
$ws.Range("C16").Value = 12.11788673706325
$ws.Range("D16").Value = 6.4714741440949
$ws.Range("E16").Value = 10.81205441551968
$ws.Range("G16").Value = 3.702858557057346
$ws.Range("I16").Value = 31.71828446208074
$ws.Range("K16").Value = 15.47416815712499
$ws.Range("M16").Value = 17.09129808124808
$ws.Range("N16").Value = 23.27740138586913

$ws.Range("B17").Value = 14.37202125942325
$ws.Range("C17").Value = 12.00363273595752
$ws.Range("D17").Value = 6.418386810633543
$ws.Range("E17").Value = 10.71275026899087
$ws.Range("G17").Value = 3.704278572873202
$ws.Range("I17").Value = 31.65147426120604
$ws.Range("K17").Value = 15.38731669001085
$ws.Range("M17").Value = 17.00939958000737
$ws.Range("N17").Value = 23.27145826298138

$ws.Range("B18").Value = 14.30859798857735
$ws.Range("C18").Value = 11.93776184308299
$ws.Range("D18").Value = 6.387834358138736
$ws.Range("E18").Value = 10.65554133300903
$ws.Range("G18").Value = 3.705106100535227
$ws.Range("I18").Value = 31.61349077953394
$ws.Range("K18").Value = 15.33757028799079
$ws.Range("M18").Value = 16.96258402555658
$ws.Range("N18").Value = 23.26820692311336

$ws.Range("B19").Value = 14.28713376521877
$ws.Range("C19").Value = 11.91543494030179
$ws.Range("D19").Value = 6.377487874417449
$ws.Range("E19").Value = 10.63615774736265
$ws.Range("G19").Value = 3.705388140502738
$ws.Range("I19").Value = 31.60070705406089
$ws.Range("K19").Value = 15.32076442935258
$ws.Range("M19").Value = 16.94678435002463
$ws.Range("N19").Value = 23.26713475400591

$ws.Range("B20").Value = 14.38376356695482
$ws.Range("C20").Value = 12.01581189736812
$ws.Range("D20").Value = 6.424040184720166
$ws.Range("E20").Value = 10.72333137533005
$ws.Range("G20").Value = 3.704126295656343
$ws.Range("I20").Value = 31.65854050883289
$ws.Range("K20").Value = 15.39654101433374
$ws.Range("M20").Value = 17.01808810585706
$ws.Range("N20").Value = 23.272073635893

$ws.Range("B21").Value = 14.70876515243896
$ws.Range("C21").Value = 12.35103424961367
$ws.Range("D21").Value = 6.580150899842055
$ws.Range("E21").Value = 11.01497059959261
$ws.Range("G21").Value = 3.700015136814427
$ws.Range("I21").Value = 31.85811393917853
$ws.Range("K21").Value = 15.65350644421695
$ws.Range("M21").Value = 17.26102265822023
$ws.Range("N21").Value = 23.29066276254567

$ws.Range("B22").Value = 14.92102933905531
$ws.Range("C22").Value = 12.5682106855194
$ws.Range("D22").Value = 6.681768813560309
$ws.Range("E22").Value = 11.20429394811011
$ws.Range("G22").Value = 3.697423466514661
$ws.Range("I22").Value = 31.99236605557427
$ws.Range("K22").Value = 15.8229347014029
$ws.Range("M22").Value = 17.4220617027899
$ws.Range("N22").Value = 23.30429875193767

$ws.Range("B23").Value = 14.80777795929783
$ws.Range("C23").Value = 12.4524990740299
$ws.Range("D23").Value = 6.627582602386335
$ws.Range("E23").Value = 11.10338720528247
$ws.Range("G23").Value = 3.698798012567103
$ws.Range("I23").Value = 31.92037040836858
$ws.Range("K23").Value = 15.73238936622759
$ws.Range("M23").Value = 17.33592012192489
$ws.Range("N23").Value = 23.29688493449041

$ws.Range("B24").Value = 14.37845481906742
$ws.Range("C24").Value = 12.01030626792858
$ws.Range("D24").Value = 6.421484390679097
$ws.Range("E24").Value = 10.71854801844752
$ws.Range("G24").Value = 3.704195105472109
$ws.Range("I24").Value = 31.65534452793234
$ws.Range("K24").Value = 15.39237011486904
$ws.Range("M24").Value = 17.01415918194661
$ws.Range("N24").Value = 23.27179491040278

$ws.Range("B25").Value = 13.91734166706539
$ws.Range("C25").Value = 11.52776179117589
$ws.Range("D25").Value = 6.198646894093115
$ws.Range("E25").Value = 10.30023474216389
$ws.Range("G25").Value = 3.710429025671207
$ws.Range("I25").Value = 31.38651264772886
$ws.Range("K25").Value = 15.03378441015925
$ws.Range("M25").Value = 16.67836751455353
$ws.Range("N25").Value = 23.25111453338247
